$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D2 to a text value "gdhf" (was numeric 6)
$ws.Range("D2").Value = "gdhf"

# Update B3 to a text value "a" (was numeric 0)
$ws.Range("B3").Value = "a"

# Update C3 to a text value "siuuuu" (was numeric 0)
$ws.Range("C3").Value = "siuuuu"
